# Implemented Pdf Data Extraction Code
# Adds Employee Number / Account Number / CIN Number / Month columns (F:I)
# to the summary report, populated for every "Attachment Found" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when Excel's normal cell-entry
# parsing would otherwise coerce it (thousands-separator numbers, mmm-yyyy
# dates, etc). Stamping the cell Text ("@") before the assignment keeps the
# literal string, then ClearFormats() drops the temporary number-format
# styling so the cell is left as a plain (unstyled) shared-string cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# New header row (F1:I1)
$ws.Range("F1").Value = "Employee Number"
$ws.Range("G1").Value = "Account Number"
$ws.Range("H1").Value = "CIN Number"
$ws.Range("I1").Value = "Month"

# Data rows: every row where an attachment/statement was found gets the new
# extracted fields. Row 3 ("Attachment Not Found") has no data, so it is
# left untouched, matching the source edit.
$dataRows = @(2, 4, 5, 6)
foreach ($r in $dataRows) {
    Set-TextValue $ws.Range("F$r") "1182"
    Set-TextValue $ws.Range("G$r") "280701501966,"
    $ws.Range("H$r").Value = "U72200TG2014PTC092878"
    Set-TextValue $ws.Range("I$r") "May-2022"
}

# Match the column widths recorded for the new columns (G/H/I only - F is
# left at its default width, matching the saved workbook). The values below
# are the closest this engine's column-width quantization (1/6-character
# steps) can land to the widths recorded in the saved file.
$ws.Columns.Item(7).ColumnWidth = 14.0
$ws.Columns.Item(8).ColumnWidth = 22.166666666666668
$ws.Columns.Item(9).ColumnWidth = 8.333333333333334

# Restore the active selection captured in the saved workbook.
$ws.Range("L13").Select()
